$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 12): update Prediction / Error values
$ws.Range("D2").Value = [double]"0.9999998921743428"
$ws.Range("E2").Value = [double]"0.9999998921743428"

# Row 3 (Control 18): update Prediction / Error values
$ws.Range("D3").Value = [double]"0.0001127654243477936"
$ws.Range("E3").Value = [double]"0.0001127654243477936"

# Row 4 (Control 34): update Prediction / Error values
$ws.Range("D4").Value = [double]"3.703921923180752E-10"
$ws.Range("E4").Value = [double]"3.703921923180752E-10"

# Row 5 (Control 42): update Prediction / Error values
$ws.Range("D5").Value = [double]"9.737364658495409E-28"
$ws.Range("E5").Value = [double]"9.737364658495409E-28"

# Row 7 (MDD 8): update Prediction / Error values
$ws.Range("D7").Value = [double]"0.0003777364295728222"
$ws.Range("E7").Value = [double]"0.9996222635704272"

# Row 9 (MDD 16): update Prediction / Error values
$ws.Range("D9").Value = [double]"0.9999605899149142"
$ws.Range("E9").Value = [double]"3.941008508578925E-05"

# Row 11 (MDD 33): update Cross Entropy Loss
$ws.Range("F11").Value = [double]"9.184396743774414"
